# Re-order the "Recorded By" (column G) author lists.
# For every row in the used range, take the comma-separated list of
# recorders stored in column G. If the first entry in that list is
# "dnasr281@gmail.com" or lower-case "system", rotate the list left by
# one position (move the first entry to the end), matching the target
# sync performed upstream. All other values (single entries, or lists
# that start with "backup@backdoor.com", "admin@admin.com" or
# capitalised "System") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if ($value -notmatch ',') { continue }

    $parts = $value -split ',\s*'
    $first = $parts[0].Trim()

    if ($first -eq 'dnasr281@gmail.com' -or $first -eq 'system') {
        $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ', '
        $cell.Value2 = $rotated
    }
}
